$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header text (Volume/Number and date range)
$ws.Range("A8").Value = "Volume 30   Number  29"
$ws.Range("C9").Value = "Report Covering the Week  7/17/2023  Through  7/23/2023"

# Update weekly crime statistics table (rows 14-30)

# Row 14
$ws.Range("C14").Value = 2
$ws.Range("D14").Value = 1
$ws.Range("E14").Value = 100
$ws.Range("I14").Value = 36
$ws.Range("J14").Value = 44
$ws.Range("K14").Value = -18.181818181818
$ws.Range("L14").Value = -26.530612244898
$ws.Range("M14").Value = -55
$ws.Range("N14").Value = -85.882352941176

# Row 15
$ws.Range("C15").Value = 2
$ws.Range("D15").Value = 6
$ws.Range("E15").Value = -66.666666666666
$ws.Range("F15").Value = 18
$ws.Range("G15").Value = 21
$ws.Range("H15").Value = -14.285714285714
$ws.Range("I15").Value = 126
$ws.Range("J15").Value = 132
$ws.Range("K15").Value = -4.545454545454
$ws.Range("L15").Value = -1.5625
$ws.Range("M15").Value = 10.526315789473
$ws.Range("N15").Value = -64.406779661017

# Row 16
$ws.Range("D16").Value = 52
$ws.Range("E16").Value = -3.846153846153
$ws.Range("F16").Value = 184
$ws.Range("G16").Value = 229
$ws.Range("H16").Value = -19.650655021834
$ws.Range("I16").Value = 1322
$ws.Range("J16").Value = 1438
$ws.Range("K16").Value = -8.066759388038
$ws.Range("L16").Value = 23.091247672253
$ws.Range("M16").Value = -30.126849894291
$ws.Range("N16").Value = -85.307846188041

# Row 17
$ws.Range("C17").Value = 105
$ws.Range("D17").Value = 91
$ws.Range("E17").Value = 15.384615384615
$ws.Range("F17").Value = 394
$ws.Range("G17").Value = 364
$ws.Range("H17").Value = 8.241758241758
$ws.Range("I17").Value = 2395
$ws.Range("J17").Value = 2291
$ws.Range("K17").Value = 4.539502400698
$ws.Range("L17").Value = 28.418230563002
$ws.Range("M17").Value = 27.597229621736
$ws.Range("N17").Value = -49.150743099787

# Row 18
$ws.Range("C18").Value = 39
$ws.Range("D18").Value = 41
$ws.Range("E18").Value = -4.878048780487
$ws.Range("F18").Value = 154
$ws.Range("G18").Value = 156
$ws.Range("H18").Value = -1.282051282051
$ws.Range("I18").Value = 1125
$ws.Range("J18").Value = 1330
$ws.Range("K18").Value = -15.413533834586
$ws.Range("L18").Value = 10.078277886497
$ws.Range("M18").Value = -30.340557275541
$ws.Range("N18").Value = -83.013739996980

# Row 19
$ws.Range("C19").Value = 105
$ws.Range("D19").Value = 120
$ws.Range("E19").Value = -12.5
$ws.Range("F19").Value = 459
$ws.Range("G19").Value = 523
$ws.Range("H19").Value = -12.237093690248
$ws.Range("I19").Value = 3186
$ws.Range("J19").Value = 3216
$ws.Range("K19").Value = -0.932835820895
$ws.Range("L19").Value = 34.487125369354
$ws.Range("M19").Value = 43.837471783295
$ws.Range("N19").Value = -12.496566877231

# Row 20
$ws.Range("C20").Value = 37
$ws.Range("D20").Value = 29
$ws.Range("E20").Value = 27.586206896551
$ws.Range("F20").Value = 159
$ws.Range("G20").Value = 151
$ws.Range("H20").Value = 5.298013245033
$ws.Range("I20").Value = 978
$ws.Range("J20").Value = 964
$ws.Range("K20").Value = 1.452282157676
$ws.Range("L20").Value = 25.063938618925
$ws.Range("M20").Value = 22.097378277153
$ws.Range("N20").Value = -81.438603150502

# Row 21
$ws.Range("C21").Value = 340
$ws.Range("D21").Value = 340
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 1372
$ws.Range("G21").Value = 1455
$ws.Range("H21").Value = -5.704467353951
$ws.Range("I21").Value = 9168
$ws.Range("J21").Value = 9415
$ws.Range("K21").Value = -2.623473181094
$ws.Range("L21").Value = 25.778570448621
$ws.Range("M21").Value = 6.679078426809
$ws.Range("N21").Value = -69.286432160804

# Row 22
$ws.Range("C22").Value = 4
$ws.Range("D22").Value = 4
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 25
$ws.Range("H22").Value = -19.354838709677
$ws.Range("I22").Value = 162
$ws.Range("J22").Value = 207
$ws.Range("K22").Value = -21.739130434782
$ws.Range("L22").Value = 25.581395348837
$ws.Range("M22").Value = -29.565217391304

# Row 23
$ws.Range("D23").Value = 31
$ws.Range("E23").Value = 3.225806451612
$ws.Range("F23").Value = 122
$ws.Range("G23").Value = 135
$ws.Range("H23").Value = -9.629629629629
$ws.Range("I23").Value = 905
$ws.Range("J23").Value = 842
$ws.Range("K23").Value = 7.482185273159
$ws.Range("L23").Value = 17.227979274611
$ws.Range("M23").Value = 45.967741935483

# Row 24
$ws.Range("C24").Value = 277
$ws.Range("D24").Value = 294
$ws.Range("E24").Value = -5.782312925170
$ws.Range("F24").Value = 1068
$ws.Range("G24").Value = 1094
$ws.Range("H24").Value = -2.376599634369
$ws.Range("I24").Value = 6937
$ws.Range("J24").Value = 7261
$ws.Range("K24").Value = -4.462195289904
$ws.Range("L24").Value = 25.443037974683
$ws.Range("M24").Value = 24.207699194270

# Row 25
$ws.Range("C25").Value = 110
$ws.Range("D25").Value = 111
$ws.Range("E25").Value = -0.900900900900
$ws.Range("F25").Value = 505
$ws.Range("G25").Value = 443
$ws.Range("H25").Value = 13.995485327313
$ws.Range("I25").Value = 3425
$ws.Range("J25").Value = 3413
$ws.Range("K25").Value = 0.351596835628
$ws.Range("L25").Value = 37.550200803212
$ws.Range("M25").Value = -23.770309370131

# Row 26
$ws.Range("C26").Value = 5
$ws.Range("E26").Value = -50
$ws.Range("F26").Value = 30
$ws.Range("G26").Value = 28
$ws.Range("H26").Value = 7.142857142857
$ws.Range("I26").Value = 195
$ws.Range("J26").Value = 207
$ws.Range("K26").Value = -5.797101449275
$ws.Range("L26").Value = -9.302325581395

# Row 27
$ws.Range("C27").Value = 18
$ws.Range("D27").Value = 9
$ws.Range("E27").Value = 100
$ws.Range("F27").Value = 51
$ws.Range("G27").Value = 58
$ws.Range("H27").Value = -12.068965517241
$ws.Range("I27").Value = 347
$ws.Range("J27").Value = 348
$ws.Range("K27").Value = -0.287356321839
$ws.Range("L27").Value = -10.567010309278

# Row 28
$ws.Range("C28").Value = 3
$ws.Range("D28").Value = 18
$ws.Range("E28").Value = -83.333333333333
$ws.Range("F28").Value = 19
$ws.Range("G28").Value = 51
$ws.Range("H28").Value = -62.745098039215
$ws.Range("I28").Value = 129
$ws.Range("J28").Value = 198
$ws.Range("K28").Value = -34.848484848484
$ws.Range("L28").Value = -39.150943396226
$ws.Range("M28").Value = -55.052264808362
$ws.Range("N28").Value = -87.875939849624

# Row 29
$ws.Range("C29").Value = 2
$ws.Range("D29").Value = 12
$ws.Range("E29").Value = -83.333333333333
$ws.Range("F29").Value = 17
$ws.Range("G29").Value = 42
$ws.Range("H29").Value = -59.523809523809
$ws.Range("I29").Value = 112
$ws.Range("J29").Value = 164
$ws.Range("K29").Value = -31.707317073170
$ws.Range("L29").Value = -37.430167597765
$ws.Range("M29").Value = -50.877192982456
$ws.Range("N29").Value = -88.308977035490

# Row 30
$ws.Range("D30").Value = 1
$ws.Range("E30").Value = 0
$ws.Range("F30").Value = 7
$ws.Range("G30").Value = 4
$ws.Range("H30").Value = 75
$ws.Range("I30").Value = 39
$ws.Range("J30").Value = 43
$ws.Range("K30").Value = -9.302325581395
$ws.Range("L30").Value = 2.631578947368
